$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "price" column (J) with header and value first so the new
# shared strings are interned before the updated date strings below.
$ws.Range("J1").Value = "price"
$ws.Range("J2").Value = "AUD `$ 250"

# Update the existing Check In / Check Out date values (stored as text)
$ws.Range("G2").Value = "29/05/2016"
$ws.Range("H2").Value = "30/05/2016"

# Give column J a sensible custom width (bestFit-style) similar to its neighbors
$ws.Columns.Item(10).ColumnWidth = 8.8

# Update the selected cell to match the saved view state
$ws.Range("H2").Select()
